$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Requisitos" entries so that LOM3231 moves below LOM3206 and LOM3215.
# Before: B26/C26=LOM3231, B27/C27=LOM3206, B28/C28=LOM3215, B29/C29=LOM3234 (unchanged)
# After : B26/C26=LOM3206, B27/C27=LOM3215, B28/C28=LOM3231, B29/C29=LOM3234 (unchanged)

$lom3206 = "LOM3206 -  Eletrônica  (Requisito)`n"
$lom3215 = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"
$lom3231 = "LOM3231 -  Métodos Experimentais da Física IV  (Indicação de Conjunto)`n"

$ws.Range("B26").Value = $lom3206
$ws.Range("C26").Value = $lom3206

$ws.Range("B27").Value = $lom3215
$ws.Range("C27").Value = $lom3215

$ws.Range("B28").Value = $lom3231
$ws.Range("C28").Value = $lom3231
